$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the driver parameters (diameter and ticks) that feed the calculation formulas
$ws.Range("C2").Value = 170
$ws.Range("C3").Value = 3200

# Update the active selection to match the edited state
$ws.Range("E7").Select()
